# Auto-generated Excel COM-interop script
# Updates profit-calculation values across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# per scheduled-runner price refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2795.5715
$ws.Range("J17").Value = 3767.25
$ws.Range("L17").Value = 11301.75
$ws.Range("N17").Value = -11637.75
$ws.Range("H19").Value = 5282
$ws.Range("I19").Value = 4783.4
$ws.Range("K19").Value = 4783.4
$ws.Range("M19").Value = -4608.4
$ws.Range("H112").Value = 5395.4
$ws.Range("I112").Value = 1512.5
$ws.Range("K112").Value = 4537.5
$ws.Range("M112").Value = -3429.5
$ws.Range("H116").Value = 5347.722
$ws.Range("I116").Value = 5632.5
$ws.Range("J116").Value = 5119.9
$ws.Range("K116").Value = 5632.5
$ws.Range("L116").Value = 5119.9
$ws.Range("M116").Value = -2190.5
$ws.Range("N116").Value = -12003.9
$ws.Range("H132").Value = 18070.441
$ws.Range("I132").Value = 20668.945
$ws.Range("K132").Value = 62006.835
$ws.Range("M132").Value = -59476.835
$ws.Range("H137").Value = 10290.289
$ws.Range("I137").Value = 20230.25
$ws.Range("J137").Value = 3061.2273
$ws.Range("K137").Value = 60690.75
$ws.Range("L137").Value = 9183.6819
$ws.Range("M137").Value = -58140.75
$ws.Range("N137").Value = -14283.6819
$ws.Range("H138").Value = 2763.7097
$ws.Range("I138").Value = 1511.375
$ws.Range("J138").Value = 4099.533
$ws.Range("K138").Value = 4534.125
$ws.Range("L138").Value = 12298.599
$ws.Range("M138").Value = 605.875
$ws.Range("N138").Value = -22578.599

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26387.354
$ws.Range("I32").Value = 27804.512
$ws.Range("K32").Value = 27804.512
$ws.Range("M32").Value = -27517.512
$ws.Range("H45").Value = 3111.524
$ws.Range("I45").Value = 2159.25
$ws.Range("K45").Value = 2159.25
$ws.Range("M45").Value = -1782.25
$ws.Range("H61").Value = 2852.2727
$ws.Range("I61").Value = 1275.037
$ws.Range("K61").Value = 1275.037
$ws.Range("M61").Value = -1063.037
$ws.Range("H74").Value = 668848.9
$ws.Range("I74").Value = 858448.5600000001
$ws.Range("K74").Value = 858448.5600000001
$ws.Range("M74").Value = -857574.5600000001
$ws.Range("H77").Value = 668848.9
$ws.Range("I77").Value = 858448.5600000001
$ws.Range("K77").Value = 4292242.800000001
$ws.Range("M77").Value = -4287874.800000001
$ws.Range("H97").Value = 1726.6571
$ws.Range("I97").Value = 1081.2916
$ws.Range("J97").Value = 3134.7273
$ws.Range("K97").Value = 1081.2916
$ws.Range("L97").Value = 3134.7273
$ws.Range("M97").Value = -585.2916
$ws.Range("N97").Value = -4126.7273
$ws.Range("H136").Value = 2852.2727
$ws.Range("I136").Value = 1275.037
$ws.Range("K136").Value = 3825.111
$ws.Range("M136").Value = -1275.111

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1819.1708
$ws.Range("I86").Value = 1416.3103
$ws.Range("K86").Value = 1416.3103
$ws.Range("M86").Value = -293.3103000000001
$ws.Range("H89").Value = 1819.1708
$ws.Range("I89").Value = 1416.3103
$ws.Range("K89").Value = 7081.5515
$ws.Range("M89").Value = -1465.5515
$ws.Range("H99").Value = 2200.8823
$ws.Range("I99").Value = 1593.2307
$ws.Range("K99").Value = 1593.2307
$ws.Range("M99").Value = -95.23070000000007
$ws.Range("H105").Value = 3491
$ws.Range("I105").Value = 3491
$ws.Range("K105").Value = 3491
$ws.Range("M105").Value = -1744
$ws.Range("H107").Value = 102754.4
$ws.Range("I107").Value = 127443
$ws.Range("K107").Value = 127443
$ws.Range("M107").Value = -125523

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3533.3333
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H34").Value = 3533.3333
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H58").Value = 15388.174
$ws.Range("I58").Value = 1340.625
$ws.Range("K58").Value = 1340.625
$ws.Range("M58").Value = -1137.625
$ws.Range("H134").Value = 2100.0557
$ws.Range("I134").Value = 1770.2258
$ws.Range("K134").Value = 5310.6774
$ws.Range("M134").Value = -2775.6774
$ws.Range("H136").Value = 15388.174
$ws.Range("I136").Value = 1340.625
$ws.Range("K136").Value = 4021.875
$ws.Range("M136").Value = -1471.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 961.44446
$ws.Range("I5").Value = 956.625
$ws.Range("K5").Value = 2869.875
$ws.Range("M5").Value = -2757.875
$ws.Range("H135").Value = 961.44446
$ws.Range("I135").Value = 956.625
$ws.Range("K135").Value = 8609.625
$ws.Range("M135").Value = -6074.625
$ws.Range("H141").Value = 4613.385
$ws.Range("J141").Value = 4000
$ws.Range("L141").Value = 12000
$ws.Range("N141").Value = -22360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6038.0557
$ws.Range("I80").Value = 3122.3333
$ws.Range("J80").Value = 11869.5
$ws.Range("K80").Value = 3122.3333
$ws.Range("L80").Value = 11869.5
$ws.Range("M80").Value = -2124.3333
$ws.Range("N80").Value = -13865.5
$ws.Range("H83").Value = 6038.0557
$ws.Range("I83").Value = 3122.3333
$ws.Range("J83").Value = 11869.5
$ws.Range("K83").Value = 15611.6665
$ws.Range("L83").Value = 59347.5
$ws.Range("M83").Value = -10619.6665
$ws.Range("N83").Value = -69331.5
$ws.Range("H102").Value = 1741.4839
$ws.Range("I102").Value = 1654.2759
$ws.Range("K102").Value = 1654.2759
$ws.Range("M102").Value = -32.27590000000009

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3928.4583
$ws.Range("I40").Value = 3716.2222
$ws.Range("K40").Value = 3716.2222
$ws.Range("M40").Value = -3580.2222
$ws.Range("H46").Value = 6475.5293
$ws.Range("I46").Value = 3267
$ws.Range("K46").Value = 3267
$ws.Range("M46").Value = -3079
$ws.Range("H55").Value = 1611.3125
$ws.Range("I55").Value = 300.5
$ws.Range("K55").Value = 300.5
$ws.Range("M55").Value = -127.5
$ws.Range("H100").Value = 5518.4
$ws.Range("I100").Value = 2624.375
$ws.Range("K100").Value = 2624.375
$ws.Range("M100").Value = -2083.375
$ws.Range("H122").Value = 8116.75
$ws.Range("I122").Value = 6245
$ws.Range("K122").Value = 18735
$ws.Range("M122").Value = -16285

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 47142.758
$ws.Range("I122").Value = 53682.6
$ws.Range("K122").Value = 161047.8
$ws.Range("M122").Value = -158597.8
